$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data to fill in for columns C (2^4), D (2^7), E (2^10), F (2^14)
# Rows: 3,4,5 (Method1 n=9,10,11), 9,10,11 (Method2 n=9,10,11), 15,16,17 (Method3 n=9,10,11)

$data = @{
    3 = @(0.000014518, 0.000386279, 0.042885777, 0.55967044)
    4 = @(0.000014518, 0.000382759, 0.008621767, 0.499698192)
    5 = @(0.000016279, 0.000383639, 0.009478796, 0.518713417)
    9 = @(0.000013639, 0.000358122, 0.007499446, 0.558310104)
    10 = @(0.000014078, 0.000359442, 0.007296188, 0.548353079)
    11 = @(0.000014079, 0.000358562, 0.007609874, 0.518421729)
    15 = @(0.000033436, 0.000288609, 0.002453621, 0.007894084)
    16 = @(0.000032557, 0.000263972, 0.001839006, 0.026222982)
    17 = @(0.000032556, 0.000272331, 0.001725497, 0.008597129)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value = $vals[0]
    $cC.NumberFormat = "0.00E+00"

    $cD = $ws.Cells.Item($row, 4)
    $cD.Value = $vals[1]
    $cD.NumberFormat = "0.00E+00"

    $cE = $ws.Cells.Item($row, 5)
    $cE.Value = $vals[2]

    $cF = $ws.Cells.Item($row, 6)
    $cF.Value = $vals[3]
}

# Update the Method 2 average formula in H12: =H18 -> =AVERAGE(H9:H11)
$ws.Range("H12").Formula = "=AVERAGE(H9:H11)"

# Column E width (bestFit, like column H)
$ws.Range("E1").EntireColumn.ColumnWidth = $ws.Range("H1").EntireColumn.ColumnWidth

# Update selection to L7
$ws.Range("L7").Select()
